$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 - Authentication implementing
$ws.Range("B26").Value = "Authentication implementing"
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 40862
$ws.Range("D25").Copy()
[void]$ws.Range("D26").PasteSpecial(-4122)

# Row 27 - User creating and deleting
$ws.Range("B27").Value = "User creating and deleting"
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 40864
$ws.Range("D25").Copy()
[void]$ws.Range("D27").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Update selection to match post-edit cursor position
[void]$ws.Range("D28").Select()
